# Edit data file column names
#
# Adds a "Parameters" header label in A1 (new shared string) and clears
# the stray number-format style override that had been left on B10 so it
# reverts to the workbook's default ("Normal") cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Parameters"
$ws.Range("B10").Style = "Normal"
